$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testcase")

# Set the source/target exclude column lists (previously blank)
$ws.Range("B13").Value = "BIRTHDATE,DEATHDATE"
$ws.Range("B26").Value = "BIRTHDATE,DEATHDATE"

# Update primary key value from "id" to "ID"
$ws.Range("B32").Value = "ID"

# Move the active cell selection from B29 to B30
$ws.Range("B30").Select()
